# Sprint 2 backlog updated
# Update the daily-hour tracking for the "Code the feature" task (row 13):
# work was logged for Tues/Wed/Thurs/Fri (columns F:I), moving each from 0 to 2
# hours. The Totals row (row 19) recalculates automatically since it holds
# SUM formulas over this range.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F13").Value = 2
$ws.Range("G13").Value = 2
$ws.Range("H13").Value = 2
$ws.Range("I13").Value = 2

# Keep the worksheet's recorded selection in sync with where editing left off.
$ws.Range("J14").Select() | Out-Null
